# "save procedure done start 진행중"
# Fill in the measurement sweep table (columns C:O, rows 1-4) on the "Sheet"
# worksheet: frequency (Hz), rf output_0 (dBm), rf output_1 (W), current (A).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet")

# Row 1: frequency sweep, 6 GHz .. 18 GHz in 1 GHz steps
$freq = @(6000000000,7000000000,8000000000,9000000000,10000000000,11000000000,12000000000,13000000000,14000000000,15000000000,16000000000,17000000000,18000000000)

# Row 2: rf output_0 sweep, 0 dBm .. 12 dBm in 1 dBm steps
$dbm = @(0,1,2,3,4,5,6,7,8,9,10,11,12)

# Row 3: rf output_1, dBm converted to Watts
$watt = @(0.001,0.001258925411794167,0.001584893192461114,0.001995262314968879,0.00251188643150958,0.003162277660168379,0.003981071705534973,0.005011872336272722,0.006309573444801934,0.007943282347242816,0.01,0.01258925411794167,0.01584893192461113)

# Row 4: current sweep, 13 .. 1 in -1 steps
$amps = @(13,12,11,10,9,8,7,6,5,4,3,2,1)

$cols = @("C","D","E","F","G","H","I","J","K","L","M","N","O")

for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "1").Value = $freq[$i]
    $ws.Range($cols[$i] + "2").Value = $dbm[$i]
    $ws.Range($cols[$i] + "3").Value = $watt[$i]
    $ws.Range($cols[$i] + "4").Value = $amps[$i]
}
